$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell
$ws.Range("A1").Value = "Datos actualizados a 16 de Septiembre de 2020 a las 19:43"

# Full refreshed country data table (País, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
# kept sorted descending by Casos totales, matching the published dashboard update
$data = @(
  ,@("Estados Unidos", 6802097, 13950, 4087090, 2514237, 0, 573, 200770)
  ,@("India", 5112431, 94397, 4017778, 1011445, 0, 1117, 83208)
  ,@("Brasil", 4392351, 8052, 3671128, 587868, 0, 148, 133355)
  ,@("Rusia", 1079519, 5670, 890114, 170488, 0, 132, 18917)
  ,@("Peru", 738020, 0, 580753, 126340, 0, 0, 30927)
  ,@("Colombia", 728590, 0, 607978, 97324, 0, 0, 23288)
  ,@("Mexico", 676487, 4771, 481068, 123741, 0, 629, 71678)
  ,@("Sudafrica", 651521, 0, 583126, 52754, 0, 0, 15641)
  ,@("España", 614360, 11193, 0, 0, 0, 239, 30243)
  ,@("Argentina", 577338, 0, 448263, 117165, 0, 58, 11910)
  ,@("Chile", 439287, 1304, 411853, 15376, 0, 18, 12058)
  ,@("Iran", 410334, 2981, 352019, 34683, 0, 179, 23632)
  ,@("Francia", 395104, 0, 89891, 274214, 0, 0, 30999)
  ,@("Reino Unido", 378219, 3991, 0, 0, 0, 20, 41684)
  ,@("Banglades", 342671, 1615, 247969, 89879, 0, 21, 4823)
  ,@("Arabia Saudita", 327551, 621, 306004, 17178, 0, 31, 4369)
  ,@("Pakistan", 303089, 665, 290760, 5936, 0, 4, 6393)
  ,@("Irak", 303059, 4357, 237241, 57570, 0, 82, 8248)
  ,@("Turquia", 296391, 1771, 262602, 26540, 0, 63, 7249)
  ,@("Italia", 291442, 1452, 215265, 40532, 0, 12, 35645)
  ,@("Filipinas", 272934, 3550, 207858, 60344, 0, 69, 4732)
  ,@("Alemania", 266541, 1697, 239100, 17994, 0, 2, 9447)
  ,@("Indonesia", 228993, 3963, 164101, 55792, 0, 135, 9100)
  ,@("Israel", 169352, 4950, 123044, 45145, 0, 16, 1163)
  ,@("Ucrania", 162660, 2958, 72324, 86996, 0, 76, 3340)
  ,@("Canada", 139421, 618, 121994, 8234, 0, 5, 9193)
  ,@("Bolivia", 128286, 667, 86216, 34623, 0, 53, 7447)
  ,@("Catar", 122449, 235, 119400, 2841, 0, 0, 208)
  ,@("Ecuador", 121525, 1972, 97063, 13466, 0, 33, 10996)
  ,@("Rumania", 107011, 1713, 43244, 59482, 0, 49, 4285)
  ,@("Kazajistan", 106984, 64, 101267, 4046, 0, 0, 1671)
  ,@("Republica Dominicana", 105521, 718, 78891, 24621, 0, 11, 2009)
  ,@("Panama", 102832, 0, 75592, 25053, 0, 0, 2187)
  ,@("Egipto", 101340, 0, 85745, 9916, 0, 0, 5679)
  ,@("Kuwait", 96999, 698, 87187, 9241, 0, 3, 571)
  ,@("Belgica", 94795, 489, 18789, 66076, 0, 3, 9930)
  ,@("Marruecos", 92016, 1692, 72968, 17362, 0, 38, 1686)
  ,@("Oman", 91196, 536, 84363, 6028, 0, 8, 805)
  ,@("Suecia", 87575, 0, 0, 0, 0, 5, 5860)
  ,@("Paises Bajos", 86320, 1542, 0, 0, 0, 2, 6260)
  ,@("China", 85214, 12, 80437, 143, 0, 0, 4634)
  ,@("Guatemala", 82684, 0, 71983, 7717, 0, 0, 2984)
  ,@("Emiratos Arabes Unidos", 81782, 842, 71456, 9924, 0, 1, 402)
  ,@("Japon", 75958, 0, 67831, 6676, 0, 0, 1451)
  ,@("Polonia", 75734, 600, 62110, 11387, 0, 10, 2237)
  ,@("Bielorrusia", 74763, 211, 72810, 1186, 0, 6, 767)
  ,@("Honduras", 68620, 0, 18487, 48046, 0, 0, 2087)
  ,@("Portugal", 65626, 605, 44528, 19220, 0, 3, 1878)
  ,@("Etiopia", 65486, 0, 25988, 38463, 0, 0, 1035)
  ,@("Venezuela", 62655, 0, 50361, 11792, 0, 0, 502)
  ,@("Barein", 61643, 0, 54831, 6596, 0, 3, 216)
  ,@("Nepal", 58327, 1539, 41706, 16242, 0, 8, 379)
  ,@("Costa Rica", 58137, 0, 21536, 35968, 0, 0, 633)
  ,@("Singapur", 57514, 26, 56884, 603, 0, 0, 27)
  ,@("Nigeria", 56478, 0, 44430, 10960, 0, 0, 1088)
  ,@("Uzbekistan", 49015, 586, 45422, 3186, 0, 5, 407)
  ,@("Argelia", 48966, 232, 34517, 12804, 0, 13, 1645)
  ,@("Suiza", 48265, 514, 39900, 6326, 0, 11, 2039)
  ,@("Armenia", 46376, 257, 42069, 3384, 0, 3, 923)
  ,@("Ghana", 45655, 0, 44797, 564, 0, 0, 294)
  ,@("Kirguistan", 45072, 73, 41210, 2799, 0, 0, 1063)
  ,@("Moldavia", 44361, 627, 32732, 10470, 0, 16, 1159)
  ,@("Chequia", 40186, 1290, 22896, 16809, 0, 5, 481)
  ,@("Afganistan", 38855, 40, 32503, 4916, 0, 10, 1436)
  ,@("Azerbaiyan", 38658, 141, 36149, 1940, 0, 3, 569)
  ,@("Kenia", 36301, 0, 23364, 12303, 0, 0, 634)
  ,@("Austria", 35073, 768, 27655, 6660, 0, 1, 758)
  ,@("Estado de Palestina", 33006, 756, 22209, 10554, 0, 14, 243)
  ,@("Serbia", 32613, 102, 31411, 466, 0, 1, 736)
  ,@("Irlanda", 31549, 0, 23364, 6398, 0, 0, 1787)
  ,@("Paraguay", 29298, 0, 15261, 13485, 0, 0, 552)
  ,@("El Salvador", 27163, 75, 19960, 6407, 0, 4, 796)
  ,@("Australia", 26779, 40, 23726, 2229, 0, 8, 824)
  ,@("Libano", 26083, 634, 9634, 16190, 0, 7, 259)
  ,@("Libia", 24936, 792, 13498, 11044, 0, 11, 394)
  ,@("Bosnia y Herzegovina", 24211, 282, 16990, 6485, 0, 11, 736)
  ,@("Corea del Sur", 22504, 113, 19310, 2827, 0, 0, 367)
  ,@("Dinamarca", 20571, 0, 16557, 3381, 0, 0, 633)
  ,@("Camerun", 20271, 0, 18837, 1019, 0, 0, 415)
  ,@("Costa de Marfil", 19100, 0, 18228, 752, 0, 0, 120)
  ,@("Bulgaria", 18216, 0, 13057, 4423, 0, 0, 736)
  ,@("Republica de Macedonia", 16088, 163, 13550, 1870, 0, 7, 668)
  ,@("Madagascar", 15871, 68, 14482, 1174, 0, 1, 215)
  ,@("Senegal", 14568, 39, 10756, 3513, 0, 1, 299)
  ,@("Hungria", 14460, 581, 4153, 9653, 0, 8, 654)
  ,@("Grecia", 14041, 311, 3804, 9921, 0, 3, 316)
  ,@("Croacia", 14029, 280, 11690, 2103, 0, 6, 236)
  ,@("Zambia", 13887, 68, 12869, 692, 0, 2, 326)
  ,@("Sudan", 13535, 0, 6759, 5940, 0, 0, 836)
  ,@("Noruega", 12431, 38, 10371, 1795, 0, 0, 265)
  ,@("Albania", 11816, 144, 6733, 4740, 0, 3, 343)
  ,@("Consejo Danes para los Refugiados", 10414, 13, 9827, 320, 0, 0, 267)
  ,@("Guinea", 10111, 0, 9444, 604, 0, 0, 63)
  ,@("Malasia", 10031, 62, 9235, 668, 0, 0, 128)
  ,@("Namibia", 9964, 63, 7502, 2354, 0, 2, 108)
  ,@("Guayana Francesa", 9578, 0, 9202, 312, 0, 0, 64)
  ,@("Maldivas", 9427, 99, 7903, 1491, 0, 0, 33)
  ,@("Tayikistan", 9171, 42, 7941, 1157, 0, 0, 73)
  ,@("Finlandia", 8750, 25, 7700, 711, 0, 0, 339)
  ,@("Gabon", 8654, 0, 7785, 816, 0, 0, 53)
  ,@("Haiti", 8530, 0, 6181, 2129, 0, 0, 220)
  ,@("Tunez", 7623, 0, 2262, 5238, 0, 0, 123)
  ,@("Zimbabue", 7576, 0, 5783, 1569, 0, 0, 224)
  ,@("Mauritania", 7319, 0, 6839, 319, 0, 0, 161)
  ,@("Luxemburgo", 7284, 0, 6561, 599, 0, 0, 124)
  ,@("Montenegro", 7061, 0, 4692, 2246, 0, 0, 123)
  ,@("Mozambique", 5994, 281, 3267, 2688, 0, 2, 39)
  ,@("Eslovaquia", 5860, 92, 3220, 2602, 0, 0, 38)
  ,@("Malaui", 5701, 0, 3762, 1761, 0, 0, 178)
  ,@("Republica de Yibuti", 5396, 0, 5331, 4, 0, 0, 61)
  ,@("Uganda", 5266, 143, 2404, 2802, 0, 2, 60)
  ,@("Suazilandia", 5128, 0, 4401, 626, 0, 0, 101)
  ,@("Guinea Ecuatorial", 5000, 0, 4496, 421, 0, 0, 83)
  ,@("Hong Kong", 4985, 9, 4663, 220, 0, 0, 102)
  ,@("Nicaragua", 4961, 0, 2913, 1901, 0, 0, 147)
  ,@("Congo", 4934, 0, 3887, 959, 0, 0, 88)
  ,@("Cabo Verde", 4904, 0, 4294, 564, 0, 0, 46)
  ,@("Cuba", 4803, 0, 4119, 576, 0, 0, 108)
  ,@("Republica de Africa Central", 4772, 0, 1828, 2882, 0, 0, 62)
  ,@("Surinam", 4625, 0, 3996, 534, 0, 0, 95)
  ,@("Ruanda", 4624, 0, 2767, 1835, 0, 0, 22)
  ,@("Jamaica", 4164, 122, 1180, 2938, 0, 0, 46)
  ,@("Eslovenia", 3954, 123, 2844, 975, 0, 0, 135)
  ,@("Jordania", 3852, 175, 2349, 1477, 0, 0, 26)
  ,@("Birmania", 3636, 134, 832, 2765, 0, 4, 39)
  ,@("Siria", 3614, 0, 871, 2583, 0, 0, 160)
  ,@("Angola", 3569, 0, 1332, 2098, 0, 0, 139)
  ,@("Tailandia", 3490, 10, 3316, 116, 0, 0, 58)
  ,@("Lituania", 3442, 45, 2125, 1230, 0, 0, 87)
  ,@("Gambia", 3428, 0, 1737, 1586, 0, 0, 105)
  ,@("Somalia", 3390, 0, 2812, 480, 0, 0, 98)
  ,@("Mayotte", 3374, 0, 2964, 370, 0, 0, 40)
  ,@("Trinidad yTobago", 3293, 70, 810, 2426, 0, 1, 57)
  ,@("Sri Lanka", 3271, 0, 3021, 237, 0, 0, 13)
  ,@("Aruba", 3152, 0, 1610, 1520, 0, 0, 22)
  ,@("Guadalupe", 3080, 0, 837, 2219, 0, 0, 24)
  ,@("Bahamas", 3032, 0, 1482, 1481, 0, 0, 69)
  ,@("Reunion", 3002, 100, 1313, 1674, 0, 0, 15)
  ,@("Mali", 2966, 26, 2311, 527, 0, 0, 128)
  ,@("Georgia", 2758, 196, 1412, 1327, 0, 0, 19)
  ,@("Estonia", 2756, 36, 2318, 374, 0, 0, 64)
  ,@("Sudan del Sur", 2594, 2, 1290, 1255, 0, 0, 49)
  ,@("Malta", 2560, 106, 1953, 591, 0, 0, 16)
  ,@("Botsuana", 2463, 0, 575, 1877, 0, 0, 11)
  ,@("Guinea-Bisau", 2275, 0, 1127, 1109, 0, 0, 39)
  ,@("Benin", 2267, 0, 1942, 285, 0, 0, 40)
  ,@("Islandia", 2189, 15, 2104, 75, 0, 0, 10)
  ,@("Sierra Leona", 2126, 0, 1640, 414, 0, 0, 72)
  ,@("Yemen", 2019, 3, 1221, 215, 0, 0, 583)
  ,@("Guyana", 1958, 0, 1302, 598, 0, 0, 58)
  ,@("Uruguay", 1827, 0, 1545, 237, 0, 0, 45)
  ,@("Nueva Zelanda", 1802, 1, 1698, 79, 0, 1, 25)
  ,@("Burkina Faso", 1733, 0, 1141, 536, 0, 0, 56)
  ,@("Togo", 1595, 0, 1219, 336, 0, 0, 40)
  ,@("Republica de Chipre", 1540, 0, 1282, 236, 0, 0, 22)
  ,@("Belice", 1528, 27, 663, 846, 0, 0, 19)
  ,@("Letonia", 1486, 4, 1248, 203, 0, 0, 35)
  ,@("Principado de Andorra", 1438, 0, 945, 440, 0, 0, 53)
  ,@("Lesoto", 1327, 0, 687, 607, 0, 0, 33)
  ,@("Liberia", 1327, 0, 1214, 31, 0, 0, 82)
  ,@("Niger", 1182, 0, 1104, 9, 0, 0, 69)
  ,@("Martinica", 1122, 0, 98, 1006, 0, 0, 18)
  ,@("Polinesia Francesa", 1099, 0, 672, 425, 0, 0, 2)
  ,@("Republica del Chad", 1087, 0, 944, 62, 0, 0, 81)
  ,@("Vietnam", 1063, 0, 936, 92, 0, 0, 35)
  ,@("Santo Tome y Principe", 907, 0, 871, 21, 0, 0, 15)
  ,@("San Marino", 723, 0, 663, 18, 0, 0, 42)
  ,@("Crucero", 712, 0, 651, 48, 0, 0, 13)
  ,@("Islas Turcas y Caicos", 650, 2, 557, 88, 0, 0, 5)
  ,@("San Martin (Parte Holandesa)", 549, 0, 468, 62, 0, 0, 19)
  ,@("Papua Nueva Guinea", 516, 5, 232, 278, 0, 0, 6)
  ,@("Tanzania", 509, 0, 183, 305, 0, 0, 21)
  ,@("Taiwan", 500, 1, 477, 16, 0, 0, 7)
  ,@("Burundi", 472, 0, 374, 97, 0, 0, 1)
  ,@("Comoras", 467, 0, 446, 14, 0, 0, 7)
  ,@("Islas Feroe", 428, 0, 412, 16, 0, 0, 0)
  ,@("Mauricio", 365, 0, 338, 17, 0, 0, 10)
  ,@("Eritrea", 364, 0, 305, 59, 0, 0, 0)
  ,@("Gibraltar", 340, 6, 309, 31, 0, 0, 0)
  ,@("Isla de Man", 339, 0, 312, 3, 0, 0, 24)
  ,@("San Martin (Parte Francesa)", 330, 0, 206, 118, 0, 0, 6)
  ,@("Mongolia", 311, 0, 301, 10, 0, 0, 0)
  ,@("Camboya", 275, 0, 274, 1, 0, 0, 0)
  ,@("Butan", 246, 0, 175, 71, 0, 0, 0)
  ,@("Islas Caimanes", 208, 0, 204, 3, 0, 0, 1)
  ,@("Curazao", 192, 23, 61, 130, 0, 0, 1)
  ,@("Barbados", 184, 0, 169, 8, 0, 0, 7)
  ,@("Monaco", 178, 0, 133, 44, 0, 0, 1)
  ,@("Bermudas", 177, 0, 164, 4, 0, 0, 9)
  ,@("Brunei", 145, 0, 139, 3, 0, 0, 3)
  ,@("Seychelles", 140, 0, 136, 4, 0, 0, 0)
  ,@("Liechtenstein", 111, 0, 105, 5, 0, 0, 1)
  ,@("Antigua y Barbuda", 95, 0, 91, 1, 0, 0, 3)
  ,@("Islas Virgenes Britanicas", 66, 0, 37, 28, 0, 0, 1)
  ,@("San Vicente y las Granadinas", 64, 0, 62, 2, 0, 0, 0)
  ,@("Macao", 46, 0, 46, 0, 0, 0, 0)
  ,@("Puerto Rico", 39, 0, 1, 36, 0, 0, 2)
  ,@("Guam", 32, 0, 0, 31, 0, 0, 1)
  ,@("Fiyi", 32, 0, 24, 6, 0, 0, 2)
  ,@("Bonaire, San Eustaquio y Saba", 28, 0, 17, 10, 0, 0, 1)
  ,@("Timor Oriental", 27, 0, 26, 1, 0, 0, 0)
  ,@("Santa Lucia", 27, 0, 26, 1, 0, 0, 0)
  ,@("Nueva Caledonia", 26, 0, 26, 0, 0, 0, 0)
  ,@("Dominica", 24, 0, 18, 6, 0, 0, 0)
  ,@("Granada", 24, 0, 24, 0, 0, 0, 0)
  ,@("San Bartolome", 23, 0, 16, 7, 0, 0, 0)
  ,@("Laos", 23, 0, 22, 1, 0, 0, 0)
  ,@("Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0)
  ,@("San Cristobal y Nieves", 17, 0, 17, 0, 0, 0, 0)
  ,@("Groenlandia", 14, 0, 14, 0, 0, 0, 0)
  ,@("Islas Malvinas", 13, 0, 13, 0, 0, 0, 0)
  ,@("Montserrat", 13, 0, 12, 0, 0, 0, 1)
  ,@("Santa Sede", 12, 0, 12, 0, 0, 0, 0)
  ,@("San Pedro y Miquelon", 11, 0, 5, 6, 0, 0, 0)
  ,@("Sahara Occidental", 10, 0, 8, 1, 0, 0, 1)
  ,@("Anguila", 3, 0, 3, 0, 0, 0, 0)
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $data[$i]
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
}
